$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 06:34"

# Row 97 - Kirguistan: refreshed stats (country unchanged)
$ws.Range("B97").Value = 906
$ws.Range("C97").Value = 11
$ws.Range("D97").Value = 650
$ws.Range("E97").Value = 244

# Row 98 - now "Consejo Danes para los Refugiados" (was "Republica de Chipre")
$ws.Range("A98").Value = "Consejo Danes para los Refugiados"
$ws.Range("B98").Value = 897
$ws.Range("C98").Value = 34
$ws.Range("D98").Value = 119
$ws.Range("E98").Value = 742
$ws.Range("F98").Value = 0
$ws.Range("H98").Value = 36

# Row 99 - now "Republica de Chipre" (was "Consejo Danes para los Refugiados")
$ws.Range("A99").Value = "Republica de Chipre"
$ws.Range("B99").Value = 889
$ws.Range("D99").Value = 400
$ws.Range("E99").Value = 474
$ws.Range("F99").Value = 15
$ws.Range("H99").Value = 15

# Row 134 - now "Congo" (was "Ruanda")
$ws.Range("A134").Value = "Congo"
$ws.Range("B134").Value = 274
$ws.Range("C134").Value = 10
$ws.Range("D134").Value = 33
$ws.Range("E134").Value = 231
$ws.Range("H134").Value = 10

# Row 135 - now "Ruanda" (was "Congo")
$ws.Range("A135").Value = "Ruanda"
$ws.Range("B135").Value = 271
$ws.Range("D135").Value = 133
$ws.Range("E135").Value = 138
$ws.Range("H135").Value = 0
